$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update res_bus/vm_pu.xlsx results: case with 380 kV done.
# Slack bus voltage setpoint changed from 1.05 to 1.02 pu, and the
# resulting per-bus voltage magnitudes were recomputed for each row
# (rows 2-25, corresponding to time steps 0-23).
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.042317241594426
$ws.Range("D2").Value = 1.044400203749006
$ws.Range("E2").Value = 1.040457699264844
$ws.Range("I2").Value = 1.040207630358473
$ws.Range("J2").Value = 1.047393714385075
$ws.Range("K2").Value = 1.047171360447787
$ws.Range("L2").Value = 1.043239992214065
$ws.Range("N2").Value = 1.04888113359724

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.043408665579052
$ws.Range("D3").Value = 1.045236430196119
$ws.Range("E3").Value = 1.041390449100146
$ws.Range("I3").Value = 1.040501075415373
$ws.Range("J3").Value = 1.048130821618611
$ws.Range("K3").Value = 1.047818751919913
$ws.Range("L3").Value = 1.04398284759735
$ws.Range("N3").Value = 1.049619287607595

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.044114784836014
$ws.Range("D4").Value = 1.045777368545556
$ws.Range("E4").Value = 1.041994239347107
$ws.Range("I4").Value = 1.040689708234716
$ws.Range("J4").Value = 1.048607110593754
$ws.Range("K4").Value = 1.048236856169141
$ws.Range("L4").Value = 1.044463139146786
$ws.Range("N4").Value = 1.050096252967714

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.044411612740209
$ws.Range("D5").Value = 1.04600474136051
$ws.Range("E5").Value = 1.042248129720167
$ws.Range("I5").Value = 1.040768711179763
$ws.Range("J5").Value = 1.048807182399437
$ws.Range("K5").Value = 1.048412435254003
$ws.Range("L5").Value = 1.044664961330264
$ws.Range("N5").Value = 1.050296608898309

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.044461450022087
$ws.Range("D6").Value = 1.046042916039021
$ws.Range("E6").Value = 1.042290762364928
$ws.Range("I6").Value = 1.040781958628255
$ws.Range("J6").Value = 1.048840765949049
$ws.Range("K6").Value = 1.048441904496392
$ws.Range("L6").Value = 1.044698842743144
$ws.Range("N6").Value = 1.050330240140414

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.044118751161126
$ws.Range("D7").Value = 1.045780406860622
$ws.Range("E7").Value = 1.041997631619159
$ws.Range("I7").Value = 1.040690765048086
$ws.Range("J7").Value = 1.048609784592005
$ws.Range("K7").Value = 1.048239203019794
$ws.Range("L7").Value = 1.044465836267785
$ws.Range("N7").Value = 1.050098930763349

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.042686115500185
$ws.Range("D8").Value = 1.044682842528707
$ws.Range("E8").Value = 1.040772876888324
$ws.Range("I8").Value = 1.04030705939746
$ws.Range("J8").Value = 1.047642961688134
$ws.Range("K8").Value = 1.047390315284671
$ws.Range("L8").Value = 1.043491123453005
$ws.Range("N8").Value = 1.049130734860058

$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.040160800393929
$ws.Range("D9").Value = 1.042747611353776
$ws.Range("E9").Value = 1.038616539246249
$ws.Range("I9").Value = 1.0396213822517
$ws.Range("J9").Value = 1.045934173955079
$ws.Range("K9").Value = 1.045888333536758
$ws.Range("L9").Value = 1.041770605095131
$ws.Range("N9").Value = 1.047419520452422

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.038476664444968
$ws.Range("D10").Value = 1.041456676811563
$ws.Range("E10").Value = 1.03718022407124
$ws.Range("I10").Value = 1.039157852459711
$ws.Range("J10").Value = 1.044791534352907
$ws.Range("K10").Value = 1.04488289194434
$ws.Range("L10").Value = 1.040621611137144
$ws.Range("N10").Value = 1.046275258170952

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.037747265545478
$ws.Range("D11").Value = 1.040897505061841
$ws.Range("E11").Value = 1.036558579493443
$ws.Range("I11").Value = 1.038955618056303
$ws.Range("J11").Value = 1.044295939314724
$ws.Range("K11").Value = 1.044446547202068
$ws.Range("L11").Value = 1.040123613536077
$ws.Range("N11").Value = 1.04577895933097

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.037476309863005
$ws.Range("D12").Value = 1.040689775528572
$ws.Range("E12").Value = 1.036327716458809
$ws.Range("I12").Value = 1.038880270310508
$ws.Range("J12").Value = 1.044111729120772
$ws.Range("K12").Value = 1.044284321706285
$ws.Range("L12").Value = 1.039938563452296
$ws.Range("N12").Value = 1.045594487537412

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.037534431897421
$ws.Range("D13").Value = 1.040734335495609
$ws.Range("E13").Value = 1.036377235398234
$ws.Range("I13").Value = 1.038896443016492
$ws.Range("J13").Value = 1.044151248453896
$ws.Range("K13").Value = 1.044319126308198
$ws.Range("L13").Value = 1.039978260565856
$ws.Range("N13").Value = 1.045634062992523

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.03772486876819
$ws.Range("D14").Value = 1.040880334649605
$ws.Range("E14").Value = 1.036539495397019
$ws.Range("I14").Value = 1.038949394457744
$ws.Range("J14").Value = 1.044280714976304
$ws.Range("K14").Value = 1.044433140604238
$ws.Range("L14").Value = 1.0401083186999
$ws.Range("N14").Value = 1.045763713372242

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.037842199942711
$ws.Range("D15").Value = 1.040970285806228
$ws.Range("E15").Value = 1.036639474892725
$ws.Range("I15").Value = 1.038981989256641
$ws.Range("J15").Value = 1.044360467117691
$ws.Range("K15").Value = 1.044503369009193
$ws.Range("L15").Value = 1.040188442316581
$ws.Range("N15").Value = 1.045843578770818

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.03852506886652
$ws.Range("D16").Value = 1.041493783233092
$ws.Range("E16").Value = 1.037221486668354
$ws.Range("I16").Value = 1.039171241997602
$ws.Range("J16").Value = 1.044824407969654
$ws.Range("K16").Value = 1.044911830010893
$ws.Range("L16").Value = 1.040654651564586
$ws.Range("N16").Value = 1.046308178472005

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.03895337186881
$ws.Range("D17").Value = 1.041822108911184
$ws.Range("E17").Value = 1.037586644728079
$ws.Range("I17").Value = 1.039289547450911
$ws.Range("J17").Value = 1.045115204876728
$ws.Range("K17").Value = 1.04516778384599
$ws.Range("L17").Value = 1.040946964867729
$ws.Range("N17").Value = 1.046599388344043

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.039203178682208
$ws.Range("D18").Value = 1.042013597524335
$ws.Range("E18").Value = 1.037799663168915
$ws.Range("I18").Value = 1.039358406083439
$ws.Range("J18").Value = 1.045284742252507
$ws.Range("K18").Value = 1.045316982656242
$ws.Range("L18").Value = 1.04111742038968
$ws.Range("N18").Value = 1.046769166482342

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.039288353818648
$ws.Range("D19").Value = 1.042078887103301
$ws.Range("E19").Value = 1.037872301689513
$ws.Range("I19").Value = 1.039381860177029
$ws.Range("J19").Value = 1.045342536632625
$ws.Range("K19").Value = 1.045367839530173
$ws.Range("L19").Value = 1.041175533539591
$ws.Range("N19").Value = 1.046827042937108

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.038907420576206
$ws.Range("D20").Value = 1.041786884554974
$ws.Range("E20").Value = 1.037547463836056
$ws.Range("I20").Value = 1.039276869594276
$ws.Range("J20").Value = 1.045084013349584
$ws.Range("K20").Value = 1.045140332221014
$ws.Range("L20").Value = 1.040915607163458
$ws.Range("N20").Value = 1.046568152521353

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.037668790538391
$ws.Range("D21").Value = 1.040837342316826
$ws.Range("E21").Value = 1.036491712672723
$ws.Range("I21").Value = 1.038933807897021
$ws.Range("J21").Value = 1.044242593740951
$ws.Range("K21").Value = 1.044399570312002
$ws.Range("L21").Value = 1.04007002179874
$ws.Range("N21").Value = 1.045725538000363

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.036889872750347
$ws.Range("D22").Value = 1.040240164027781
$ws.Range("E22").Value = 1.035828171327339
$ws.Range("I22").Value = 1.038716787101816
$ws.Range("J22").Value = 1.043712841532071
$ws.Range("K22").Value = 1.043932969700143
$ws.Range("L22").Value = 1.039537954383639
$ws.Range("N22").Value = 1.045195033482583

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.037302805553523
$ws.Range("D23").Value = 1.04055675499896
$ws.Range("E23").Value = 1.036179903339603
$ws.Range("I23").Value = 1.038831959485076
$ws.Range("J23").Value = 1.04399374137239
$ws.Range("K23").Value = 1.04418040451853
$ws.Range("L23").Value = 1.039820052709636
$ws.Range("N23").Value = 1.045476332232894

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.03892818403367
$ws.Range("D24").Value = 1.041802800981454
$ws.Range("E24").Value = 1.037565167908227
$ws.Range("I24").Value = 1.039282598625905
$ws.Range("J24").Value = 1.045098107703674
$ws.Range("K24").Value = 1.045152736721269
$ws.Range("L24").Value = 1.040929776501866
$ws.Range("N24").Value = 1.046582266891042

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.040813756695373
$ws.Range("D25").Value = 1.043248053618289
$ws.Range("E25").Value = 1.039173785685398
$ws.Range("I25").Value = 1.039799776454039
$ws.Range("J25").Value = 1.046376543948611
$ws.Range("K25").Value = 1.046277358296027
$ws.Range("L25").Value = 1.042215749731474
$ws.Range("N25").Value = 1.047862518662085

Write-Host "done"